$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1156546.8
$ws.Range("J17").Value = 1156546.8
$ws.Range("L17").Value = 3469640.4
$ws.Range("N17").Value = -3469976.4
$ws.Range("H32").Value = 1233.5
$ws.Range("I32").Value = 401
$ws.Range("J32").Value = 1400
$ws.Range("K32").Value = 401
$ws.Range("L32").Value = 1400
$ws.Range("M32").Value = -75
$ws.Range("N32").Value = -2052
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H62").Value = 4505
$ws.Range("I62").Value = 3257.5
$ws.Range("J62").Value = 7000
$ws.Range("K62").Value = 3257.5
$ws.Range("L62").Value = 7000
$ws.Range("M62").Value = -2633.5
$ws.Range("N62").Value = -8248
$ws.Range("H64").Value = 9000
$ws.Range("J64").Value = 9000
$ws.Range("L64").Value = 9000
$ws.Range("N64").Value = -9496
$ws.Range("H65").Value = 4505
$ws.Range("I65").Value = 3257.5
$ws.Range("J65").Value = 7000
$ws.Range("K65").Value = 16287.5
$ws.Range("L65").Value = 35000
$ws.Range("M65").Value = -13167.5
$ws.Range("N65").Value = -41240
$ws.Range("H67").Value = 9000
$ws.Range("J67").Value = 9000
$ws.Range("L67").Value = 9000
$ws.Range("N67").Value = -10716
$ws.Range("H98").Value = 2076.9355
$ws.Range("I98").Value = 2149.7307
$ws.Range("K98").Value = 2149.7307
$ws.Range("M98").Value = -651.7307000000001
$ws.Range("H100").Value = 2895.6785
$ws.Range("I100").Value = 1840.5883
$ws.Range("K100").Value = 1840.5883
$ws.Range("M100").Value = -1299.5883
$ws.Range("H111").Value = 1333.5555
$ws.Range("J111").Value = 1810.6
$ws.Range("L111").Value = 5431.799999999999
$ws.Range("N111").Value = -11565.8
$ws.Range("H113").Value = 21894.625
$ws.Range("I113").Value = 39990.25
$ws.Range("K113").Value = 39990.25
$ws.Range("M113").Value = -36736.25
$ws.Range("H122").Value = 2076.9355
$ws.Range("I122").Value = 2149.7307
$ws.Range("K122").Value = 6449.1921
$ws.Range("M122").Value = -3999.1921
$ws.Range("H131").Value = 557347
$ws.Range("I131").Value = 626016
$ws.Range("K131").Value = 1878048
$ws.Range("M131").Value = -1873008
$ws.Range("H135").Value = 1204.9584
$ws.Range("I135").Value = 1239.0476
$ws.Range("J135").Value = 966.3333
$ws.Range("K135").Value = 11151.4284
$ws.Range("L135").Value = 8696.9997
$ws.Range("M135").Value = -8616.428400000001
$ws.Range("N135").Value = -13766.9997

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 1186.6666
$ws.Range("J12").Value = 550
$ws.Range("L12").Value = 550
$ws.Range("N12").Value = -896
$ws.Range("H19").Value = 8354
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H29").Value = 14584.333
$ws.Range("I29").Value = 27689.75
$ws.Range("J29").Value = 4100
$ws.Range("K29").Value = 27689.75
$ws.Range("L29").Value = 4100
$ws.Range("M29").Value = -27381.75
$ws.Range("N29").Value = -4716
$ws.Range("H32").Value = 3790254.5
$ws.Range("I32").Value = 4546760.5
$ws.Range("K32").Value = 4546760.5
$ws.Range("M32").Value = -4546473.5
$ws.Range("H43").Value = 23698.777
$ws.Range("I43").Value = 18848
$ws.Range("J43").Value = 24305.125
$ws.Range("K43").Value = 18848
$ws.Range("L43").Value = 24305.125
$ws.Range("M43").Value = -18535
$ws.Range("N43").Value = -24931.125
$ws.Range("H45").Value = 1763.8889
$ws.Range("I45").Value = 1979.3334
$ws.Range("K45").Value = 1979.3334
$ws.Range("M45").Value = -1602.3334
$ws.Range("H103").Value = 44750
$ws.Range("J103").Value = 44750
$ws.Range("L103").Value = 44750
$ws.Range("N103").Value = -47094
$ws.Range("H122").Value = 1530
$ws.Range("I122").Value = 1569.6666
$ws.Range("K122").Value = 4708.9998
$ws.Range("M122").Value = -2258.9998
$ws.Range("H132").Value = 4911416.5
$ws.Range("I132").Value = 2281.3845
$ws.Range("J132").Value = 9469899
$ws.Range("K132").Value = 6844.1535
$ws.Range("L132").Value = 28409697
$ws.Range("M132").Value = -4314.1535
$ws.Range("N132").Value = -28414757

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 956
$ws.Range("I64").Value = 1000
$ws.Range("J64").Value = 945
$ws.Range("K64").Value = 1000
$ws.Range("L64").Value = 945
$ws.Range("M64").Value = -775
$ws.Range("N64").Value = -1395
$ws.Range("H67").Value = 956
$ws.Range("I67").Value = 1000
$ws.Range("J67").Value = 945
$ws.Range("K67").Value = 1000
$ws.Range("L67").Value = 945
$ws.Range("M67").Value = -220
$ws.Range("N67").Value = -2505
$ws.Range("H82").Value = 9866.5
$ws.Range("I82").Value = 2839.8
$ws.Range("K82").Value = 2839.8
$ws.Range("M82").Value = -2456.8
$ws.Range("H85").Value = 9866.5
$ws.Range("I85").Value = 2839.8
$ws.Range("K85").Value = 2839.8
$ws.Range("M85").Value = -1513.8
$ws.Range("H86").Value = 1742.7646
$ws.Range("I86").Value = 1524.2307
$ws.Range("J86").Value = 2453
$ws.Range("K86").Value = 1524.2307
$ws.Range("L86").Value = 2453
$ws.Range("M86").Value = -401.2307000000001
$ws.Range("N86").Value = -4699
$ws.Range("H89").Value = 1742.7646
$ws.Range("I89").Value = 1524.2307
$ws.Range("J89").Value = 2453
$ws.Range("K89").Value = 7621.1535
$ws.Range("L89").Value = 12265
$ws.Range("M89").Value = -2005.1535
$ws.Range("N89").Value = -23497
$ws.Range("H94").Value = 1220.9667
$ws.Range("I94").Value = 1339.762
$ws.Range("J94").Value = 943.7778
$ws.Range("K94").Value = 1339.762
$ws.Range("L94").Value = 943.7778
$ws.Range("M94").Value = -888.7619999999999
$ws.Range("N94").Value = -1845.7778
$ws.Range("H100").Value = 23167.428
$ws.Range("J100").Value = 23167.428
$ws.Range("L100").Value = 23167.428
$ws.Range("N100").Value = -25331.428
$ws.Range("H102").Value = 14350.125
$ws.Range("I102").Value = 14350.125
$ws.Range("K102").Value = 14350.125
$ws.Range("M102").Value = -11105.125
$ws.Range("H138").Value = 67278.57000000001
$ws.Range("J138").Value = 70146.16
$ws.Range("L138").Value = 70146.16
$ws.Range("N138").Value = -80426.16

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 815.6
$ws.Range("J7").Value = 389.66666
$ws.Range("L7").Value = 389.66666
$ws.Range("N7").Value = -615.66666
$ws.Range("H22").Value = 959.63635
$ws.Range("I22").Value = 457.57144
$ws.Range("J22").Value = 1838.25
$ws.Range("K22").Value = 457.57144
$ws.Range("L22").Value = 1838.25
$ws.Range("M22").Value = -107.57144
$ws.Range("N22").Value = -2538.25
$ws.Range("H31").Value = 10256.091
$ws.Range("I31").Value = 960.3461
$ws.Range("J31").Value = 44783.145
$ws.Range("K31").Value = 960.3461
$ws.Range("L31").Value = 44783.145
$ws.Range("M31").Value = -665.3461
$ws.Range("N31").Value = -45373.145
$ws.Range("H34").Value = 10256.091
$ws.Range("I34").Value = 960.3461
$ws.Range("J34").Value = 44783.145
$ws.Range("K34").Value = 960.3461
$ws.Range("L34").Value = 44783.145
$ws.Range("M34").Value = -758.3461
$ws.Range("N34").Value = -45187.145
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("L44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("N44").ClearContents()
$ws.Range("H58").Value = 16359.259
$ws.Range("I58").Value = 8926.75
$ws.Range("J58").Value = 22305.268
$ws.Range("K58").Value = 8926.75
$ws.Range("L58").Value = 22305.268
$ws.Range("M58").Value = -8723.75
$ws.Range("N58").Value = -22711.268
$ws.Range("H87").Value = 46776.332
$ws.Range("J87").Value = 46776.332
$ws.Range("L87").Value = 46776.332
$ws.Range("N87").Value = -49148.332
$ws.Range("H90").Value = 46776.332
$ws.Range("J90").Value = 46776.332
$ws.Range("L90").Value = 140328.996
$ws.Range("N90").Value = -152184.996
$ws.Range("H122").Value = 4072.7144
$ws.Range("I122").Value = 2500
$ws.Range("K122").Value = 7500
$ws.Range("M122").Value = -5050
$ws.Range("H132").Value = 26399942
$ws.Range("I132").Value = 2399.9062
$ws.Range("K132").Value = 7199.7186
$ws.Range("M132").Value = -4669.7186
$ws.Range("H136").Value = 16359.259
$ws.Range("I136").Value = 8926.75
$ws.Range("J136").Value = 22305.268
$ws.Range("K136").Value = 26780.25
$ws.Range("L136").Value = 66915.804
$ws.Range("M136").Value = -24230.25
$ws.Range("N136").Value = -72015.804

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H26").Value = 1950.6364
$ws.Range("I26").Value = 2509
$ws.Range("K26").Value = 7527
$ws.Range("M26").Value = -7239
$ws.Range("H38").Value = 155.1875
$ws.Range("I38").Value = 154.125
$ws.Range("J38").Value = 156.25
$ws.Range("K38").Value = 462.375
$ws.Range("L38").Value = 468.75
$ws.Range("M38").Value = -115.375
$ws.Range("N38").Value = -1162.75
$ws.Range("H81").Value = 8502502
$ws.Range("J81").Value = 10202003
$ws.Range("L81").Value = 30606009
$ws.Range("N81").Value = -30608255
$ws.Range("H84").Value = 8502502
$ws.Range("J84").Value = 10202003
$ws.Range("L84").Value = 91818027
$ws.Range("N84").Value = -91829259
$ws.Range("H103").Value = 1780
$ws.Range("J103").Value = 1725
$ws.Range("L103").Value = 5175
$ws.Range("N103").Value = -6933
$ws.Range("H115").Value = 2790.6
$ws.Range("I115").Value = 1984.3334
$ws.Range("J115").Value = 4000
$ws.Range("K115").Value = 5953.0002
$ws.Range("L115").Value = 12000
$ws.Range("M115").Value = -4778.0002
$ws.Range("N115").Value = -14350
$ws.Range("H131").Value = 1479.72
$ws.Range("I131").Value = 0
$ws.Range("J131").Value = 1479.72
$ws.Range("K131").Value = 0
$ws.Range("L131").Value = 4439.16
$ws.Range("M131").ClearContents()
$ws.Range("N131").Value = -14519.16

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 3233.3333
$ws.Range("J6").Value = 3233.3333
$ws.Range("L6").Value = 3233.3333
$ws.Range("N6").Value = -3459.3333
$ws.Range("H16").Value = 3233.3333
$ws.Range("J16").Value = 3233.3333
$ws.Range("L16").Value = 3233.3333
$ws.Range("N16").Value = -3733.3333
$ws.Range("H21").Value = 18750
$ws.Range("I21").Value = 15000
$ws.Range("K21").Value = 15000
$ws.Range("M21").Value = -14827
$ws.Range("H30").Value = 18750
$ws.Range("I30").Value = 15000
$ws.Range("K30").Value = 15000
$ws.Range("M30").Value = -14895
$ws.Range("H122").Value = 3051.923
$ws.Range("I122").Value = 3999.6667
$ws.Range("J122").Value = 2767.6
$ws.Range("K122").Value = 11999.0001
$ws.Range("L122").Value = 8302.799999999999
$ws.Range("M122").Value = -9549.000100000001
$ws.Range("N122").Value = -13202.8
$ws.Range("H123").Value = 56058.8
$ws.Range("J123").Value = 56058.8
$ws.Range("L123").Value = 56058.8
$ws.Range("N123").Value = -60958.8
$ws.Range("H126").Value = 8248.842000000001
$ws.Range("I126").Value = 9339.866
$ws.Range("K126").Value = 28019.598
$ws.Range("M126").Value = -25549.598

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 7799.1333
$ws.Range("I13").Value = 5000
$ws.Range("J13").Value = 7999.0713
$ws.Range("K13").Value = 5000
$ws.Range("L13").Value = 7999.0713
$ws.Range("M13").Value = -4860
$ws.Range("N13").Value = -8279.0713
$ws.Range("H26").Value = 46599.8
$ws.Range("I26").Value = 45749.75
$ws.Range("K26").Value = 45749.75
$ws.Range("M26").Value = -45454.75
$ws.Range("H40").Value = 1461.6364
$ws.Range("I40").Value = 1407.8
$ws.Range("J40").Value = 2000
$ws.Range("K40").Value = 1407.8
$ws.Range("L40").Value = 2000
$ws.Range("M40").Value = -1271.8
$ws.Range("N40").Value = -2272
$ws.Range("H46").Value = 2607.75
$ws.Range("J46").Value = 2979.3
$ws.Range("L46").Value = 2979.3
$ws.Range("N46").Value = -3355.3
$ws.Range("H61").Value = 4114
$ws.Range("I61").Value = 3759.6
$ws.Range("K61").Value = 3759.6
$ws.Range("M61").Value = -3557.6
$ws.Range("H113").Value = 4114
$ws.Range("I113").Value = 3759.6
$ws.Range("K113").Value = 3759.6
$ws.Range("M113").Value = -1589.6
$ws.Range("H122").Value = 6205
$ws.Range("I122").Value = 6560
$ws.Range("K122").Value = 19680
$ws.Range("M122").Value = -17230
$ws.Range("H132").Value = 1944054.8
$ws.Range("I132").Value = 1897.7916
$ws.Range("K132").Value = 5693.3748
$ws.Range("M132").Value = -3163.3748

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 5000
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -4480
$ws.Range("N54").ClearContents()
$ws.Range("H81").Value = 530.3333
$ws.Range("I81").Value = 536.4
$ws.Range("J81").Value = 500
$ws.Range("K81").Value = 1072.8
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = -11.79999999999995
$ws.Range("N81").Value = -3122
$ws.Range("H84").Value = 530.3333
$ws.Range("I84").Value = 536.4
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 5364
$ws.Range("L84").Value = 5000
$ws.Range("M84").Value = -60
$ws.Range("N84").Value = -15608
$ws.Range("H86").Value = 3391666
$ws.Range("J86").Value = 3391666
$ws.Range("L86").Value = 3391666
$ws.Range("N86").Value = -3393912
$ws.Range("H89").Value = 3391666
$ws.Range("J89").Value = 3391666
$ws.Range("L89").Value = 16958330
$ws.Range("N89").Value = -16969562
$ws.Range("H107").Value = 2313.8333
$ws.Range("I107").Value = 2313.8333
$ws.Range("K107").Value = 6941.499899999999
$ws.Range("M107").Value = -5021.499899999999
$ws.Range("H122").Value = 2703.7778
$ws.Range("I122").Value = 2784.7144
$ws.Range("K122").Value = 8354.143199999999
$ws.Range("M122").Value = -5904.143199999999
$ws.Range("H132").Value = 550256.4
$ws.Range("I132").Value = 2072.6875
$ws.Range("J132").Value = 2742991
$ws.Range("K132").Value = 6218.0625
$ws.Range("L132").Value = 8228973
$ws.Range("M132").Value = -3688.0625
$ws.Range("N132").Value = -8234033
$ws.Range("H136").Value = 340510.28
$ws.Range("I136").Value = 2010.9546
$ws.Range("J136").Value = 1167953.1
$ws.Range("K136").Value = 6032.8638
$ws.Range("L136").Value = 3503859.3
$ws.Range("M136").Value = -3482.8638
$ws.Range("N136").Value = -3508959.3
